$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 400, shifting existing rows 400-456 down to 401-457
$ws.Rows("400:400").Insert()

# Populate new row 400 with the new record's data
$ws.Range("A400").Value = 9
$ws.Range("B400").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C400").Value = "Metropolitana"
$ws.Range("D400").Value = 45218
$ws.Range("D400").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E400").Value = 13
$ws.Range("F400").Value = 100112021
$ws.Range("G400").Value = "Ají"
$ws.Range("H400").Value = "Inferno"
$ws.Range("I400").Value = "Primera"
$ws.Range("J400").Value = 70
$ws.Range("K400").Value = 25000
$ws.Range("L400").Value = 28000
$ws.Range("M400").Value = 26500
$ws.Range("N400").Value = "`$/caja 10 kilos"
$ws.Range("O400").Value = "Región de Arica y Parinacota"
$ws.Range("P400").Value = 2650
$ws.Range("Q400").Value = 10
$ws.Range("R400").Value = "Hortaliza"
